# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.076.00"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'1.835.09"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'244.60"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "'0.6336"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.07539"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").Value = "'0.2945"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "'22.93"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "'0.07740"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "'1.834.95"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'5.005"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "'83.24"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'0.000009602"
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("D17").Value = "'6.095"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "'29.096.38"
$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "'226.38"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'7.200"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'160.67"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("D26").Value = "'8.551"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "'17.94"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("D30").Value = "'4.073"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'0.05384"
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("D33").Value = "'1.865"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "'0.7452"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").Value = "'2.660"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "'1.242.73"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'2.758"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").Value = "'0.9069"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'101.98"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'1.985.98"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").Value = "'64.86"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").Value = "'0.5116"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'0.4094"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").Value = "'9.073"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").Value = "'6.777"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("E51").Value = "  +0.44%  "
